# Fall 2020 schedule table update
#
# - drop the "Facilitators" column (4th column) and widen the
#   remaining "Topic" / "Presenter" columns to take up the freed space
# - turn the "Week of Oct. 5" placeholder row into a real "October 9, 2020"
#   entry and tidy up the garbled presenter name
# - fill in the four still-empty rows with the rest of the Fall semester
#   schedule (Oct 23, Nov 6, Nov 20, Dec 4)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellXml {
    param($Row, $Col, $InnerXml)

    $cell = $t.Cell($Row, $Col)
    $range = $cell.Range
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + $InnerXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$range.InsertXML($pkg)
}

function Run-Text {
    param($Text)
    return '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">' + $Text + '</w:t></w:r>'
}

# ---------------------------------------------------------------------
# 1. Drop the 4th ("Facilitators") column and resize the others so the
#    grid goes from 2284/2735/2322/2009 dxa to 2284/4515/2410 dxa.
# ---------------------------------------------------------------------
$t.Columns.Item(4).Delete()
$t.Columns.Item(2).Width = 4515 / 20.0
$t.Columns.Item(3).Width = 2410 / 20.0

# ---------------------------------------------------------------------
# 2. Row 2: "Week of Oct. 5" -> "October 9, 2020"; clean up "Gwen?" ->
#    "Gwen". The date cell is rewritten wholesale (dropping the stray
#    "_GoBack" bookmark that used to sit after the date text - it is
#    re-created further down, on the November 6 row); the presenter
#    cell is a plain text-only fix done with Find/Replace so the
#    existing paragraph/run metadata there is left untouched.
# ---------------------------------------------------------------------
$row2Tc1 = Run-Text "October 9, 2020"
Set-CellXml 2 1 $row2Tc1
[void]$d.Content.Find.Execute("Gwen?", $true, $false, $false, $false, $false, $true, 1, $false, "Gwen", 2)

# ---------------------------------------------------------------------
# 3. Row 3 (was blank): October 23, 2020 - Project Workflow - Elizabeth
#    Morin-Lessard.
# ---------------------------------------------------------------------
$row3Tc1 = Run-Text "October 23, 2020"
$row3Tc2 = (Run-Text "Project ") + (Run-Text "W") + (Run-Text "orkflow")
$row3Tc3 = Run-Text "Elizabeth Morin-Lessard"
Set-CellXml 3 1 $row3Tc1
Set-CellXml 3 2 $row3Tc2
Set-CellXml 3 3 $row3Tc3

# ---------------------------------------------------------------------
# 4. Row 4 (was blank): November 6, 2020 - Open Data, Data Sharing and
#    Ethics - Jenelle Morgan & Chelsea Moran. The "_GoBack" bookmark
#    (Word's last-edit marker) now lives between "6" and ", 2020".
# ---------------------------------------------------------------------
$row4Tc1 = (Run-Text "November ") + (Run-Text "6") +
           '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
           (Run-Text ", 2020")
$row4Tc2 = (Run-Text "Open Data") + (Run-Text ", Data Sharing") + (Run-Text " and Ethics")
$row4Tc3 = Run-Text "Jenelle Morgan &amp; Chelsea Moran"
Set-CellXml 4 1 $row4Tc1
Set-CellXml 4 2 $row4Tc2
Set-CellXml 4 3 $row4Tc3

# ---------------------------------------------------------------------
# 5. Row 5 (was blank): November 20, 2020 - Pre-registration - Emiko
#    Muraki & Brittany Lindsay.
# ---------------------------------------------------------------------
$row5Tc1 = Run-Text "November 20, 2020"
$row5Tc2 = Run-Text "Pre-registration"
$row5Tc3 = (Run-Text "Emiko ") +
           '<w:proofErr w:type="spellStart"/>' + (Run-Text "Muraki") + '<w:proofErr w:type="spellEnd"/>' +
           (Run-Text " &amp; Brittany Lindsay")
Set-CellXml 5 1 $row5Tc1
Set-CellXml 5 2 $row5Tc2
Set-CellXml 5 3 $row5Tc3

# ---------------------------------------------------------------------
# 6. Row 6 (was blank): December 4, 2020 - Troubleshooting Open Science
#    - Rosemary Twomey.
# ---------------------------------------------------------------------
$row6Tc1 = Run-Text "December 4, 2020"
$row6Tc2 = Run-Text "Troubleshooting Open Science"
$row6Tc3 = Run-Text "Rosemary Twomey"
Set-CellXml 6 1 $row6Tc1
Set-CellXml 6 2 $row6Tc2
Set-CellXml 6 3 $row6Tc3
